# question11.xlsx -- "finally done with 11th question"
#
# 1) Sheet "Initial Conditions": fix the G-column slope formulas so they
#    divide by the *interval width* (F_i - F_(i-1)) instead of the
#    cumulative-frequency column E, and add a new note row (merged,
#    "Neutral" cell style) below the table.
# 2) Sheet "Simulation Result": fill in the generated-solution numbers for
#    random-number columns D..H (rows 6-10), which were previously blank.
# 3) Restore each sheet's active-cell selection to where the author left it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Initial Conditions")
$ws2 = $wb.Worksheets.Item("Simulation Result")

# ---------------------------------------------------------------------
# 1) Slope formulas in column G (the "a_i" column) on the first sheet.
#    G6 has no previous interval, so it still divides by F6 alone; every
#    row below divides by the gap between consecutive cumulative values.
# ---------------------------------------------------------------------
$ws1.Range("G6").Formula     = "= 15/F6"
$ws1.Range("G7").Formula     = "= 15/(F7-F6)"
$ws1.Range("G8:G12").Formula = "= 15/(F8-F7)"

# ---------------------------------------------------------------------
# New note below the table, styled with the built-in "Neutral" cell
# style and merged across C19:F19, same footprint as the "Good"-styled
# note already sitting in C18:F18.
# ---------------------------------------------------------------------
$ws1.Range("C19").Value = "And solutions for the selceted random numbers in next sheet"
$ws1.Range("C19:F19").Style = "Neutral"
$ws1.Range("C19:F19").Merge()

# ---------------------------------------------------------------------
# 2) Fill in the rest of the solved-for table on "Simulation Result"
#    (columns D-H, rows 6-10) -- column C was already populated.
# ---------------------------------------------------------------------
$solved = @(
    @(48.001100000000001, 47.400300000000001, 47.220100000000002, 47.165999999999997, 47.149799999999999),
    @(156.97999999999999, 156.36799999999999, 156.13800000000001, 156.05199999999999, 156.01900000000001),
    @(200.68799999999999, 193.39099999999999, 192.78299999999999, 192.732,             192.72800000000001),
    @(197.55699999999999, 193.13,              192.761,            192.73,             192.72800000000001),
    @(106.69499999999999, 104.17400000000001, 103.54300000000001, 103.386,             103.346)
)
for ($r = 0; $r -lt $solved.Length; $r++) {
    $row = $solved[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item(6 + $r, 4 + $c).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# 3) Selections: author ended up on G23 on sheet 1 and G18 on sheet 2,
#    with "Initial Conditions" left as the active tab.
# ---------------------------------------------------------------------
$ws2.Range("G18").Select()
$ws1.Range("G23").Select()
